$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51. This shifts the existing rows 51-68
# down to 52-69, preserving their contents/formatting, and extends the
# sheet dimension to A1:R69 automatically.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new data record.
$ws.Range("A51").Value = 5
$ws.Range("B51").Value = "Macroferia Regional de Talca"
$ws.Range("C51").Value = "Maule"
$ws.Range("D51").Value = 44839
$ws.Range("E51").Value = 7
$ws.Range("F51").Value = 300000000
$ws.Range("G51").Value = "Espárragos"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 1300
$ws.Range("L51").Value = 1300
$ws.Range("M51").Value = 1300
$ws.Range("N51").Value = "$/kilo"
$ws.Range("O51").Value = "Provincia de Linares"
$ws.Range("P51").Value = 1300
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"
